$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1765.3846
$ws.Range("I15").Value = 1765.3846
$ws.Range("K15").Value = 5296.1538
$ws.Range("M15").Value = -5127.1538
$ws.Range("H19").Value = 1092.1428
$ws.Range("J19").Value = 1032.6666
$ws.Range("L19").Value = 1032.6666
$ws.Range("N19").Value = -1382.6666
$ws.Range("H28").Value = 6921.533
$ws.Range("I28").Value = 7935.3076
$ws.Range("J28").Value = 332
$ws.Range("K28").Value = 7935.3076
$ws.Range("L28").Value = 332
$ws.Range("M28").Value = -7450.3076
$ws.Range("N28").Value = -1302
$ws.Range("H43").Value = 5917.8076
$ws.Range("I43").Value = 5055.3335
$ws.Range("J43").Value = 6374.4116
$ws.Range("K43").Value = 5055.3335
$ws.Range("L43").Value = 6374.4116
$ws.Range("M43").Value = -4986.3335
$ws.Range("N43").Value = -6512.4116
$ws.Range("H100").Value = 4872.8887
$ws.Range("I100").Value = 2171.3
$ws.Range("J100").Value = 8249.875
$ws.Range("K100").Value = 2171.3
$ws.Range("L100").Value = 8249.875
$ws.Range("M100").Value = -1630.3
$ws.Range("N100").Value = -9331.875
$ws.Range("H113").Value = 3555.4285
$ws.Range("I113").Value = 3322
$ws.Range("J113").Value = 4255.7144
$ws.Range("K113").Value = 3322
$ws.Range("L113").Value = 4255.7144
$ws.Range("M113").Value = -68
$ws.Range("N113").Value = -10763.7144
$ws.Range("H137").Value = 2342.625
$ws.Range("I137").Value = 2258
$ws.Range("K137").Value = 6774
$ws.Range("M137").Value = -4224

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2485.2
$ws.Range("I2").Value = 1538.1
$ws.Range("K2").Value = 1538.1
$ws.Range("M2").Value = -1425.1
$ws.Range("H32").Value = 5376.8125
$ws.Range("I32").Value = 4812.2334
$ws.Range("K32").Value = 4812.2334
$ws.Range("M32").Value = -4525.2334
$ws.Range("H50").Value = 5634.8
$ws.Range("I50").Value = 449.33334
$ws.Range("J50").Value = 7857.143
$ws.Range("K50").Value = 449.33334
$ws.Range("L50").Value = 7857.143
$ws.Range("M50").Value = 264.66666
$ws.Range("N50").Value = -9285.143
$ws.Range("H61").Value = 1255.4
$ws.Range("I61").Value = 1233.5
$ws.Range("K61").Value = 1233.5
$ws.Range("M61").Value = -1021.5
$ws.Range("H97").Value = 358.6111
$ws.Range("I97").Value = 204.61539
$ws.Range("J97").Value = 759
$ws.Range("K97").Value = 204.61539
$ws.Range("L97").Value = 759
$ws.Range("M97").Value = 291.38461
$ws.Range("N97").Value = -1751
$ws.Range("H116").Value = 2485.2
$ws.Range("I116").Value = 1538.1
$ws.Range("K116").Value = 1538.1
$ws.Range("M116").Value = 755.9000000000001
$ws.Range("H123").Value = 112666.664
$ws.Range("J123").Value = 60000
$ws.Range("L123").Value = 60000
$ws.Range("N123").Value = -69800
$ws.Range("H136").Value = 1255.4
$ws.Range("I136").Value = 1233.5
$ws.Range("K136").Value = 3700.5
$ws.Range("M136").Value = -1150.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2485.2
$ws.Range("I3").Value = 1538.1
$ws.Range("K3").Value = 1538.1
$ws.Range("M3").Value = -1424.1
$ws.Range("H22").Value = 1046.25
$ws.Range("I22").Value = 1046.25
$ws.Range("K22").Value = 1046.25
$ws.Range("M22").Value = -873.25
$ws.Range("H86").Value = 9971.25
$ws.Range("J86").Value = 9995
$ws.Range("L86").Value = 9995
$ws.Range("N86").Value = -12241
$ws.Range("H89").Value = 9971.25
$ws.Range("J89").Value = 9995
$ws.Range("L89").Value = 49975
$ws.Range("N89").Value = -61207
$ws.Range("H94").Value = 4400
$ws.Range("H99").Value = 2990.45
$ws.Range("I99").Value = 1723.2222
$ws.Range("K99").Value = 1723.2222
$ws.Range("M99").Value = -225.2221999999999
$ws.Range("H134").Value = 1488.3549
$ws.Range("I134").Value = 1147.4584
$ws.Range("K134").Value = 3442.3752
$ws.Range("M134").Value = -907.3751999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 1705.5555
$ws.Range("I10").Value = 1700
$ws.Range("J10").Value = 1725
$ws.Range("K10").Value = 1700
$ws.Range("L10").Value = 1725
$ws.Range("M10").Value = -1561
$ws.Range("N10").Value = -2003
$ws.Range("H21").Value = 5006.5
$ws.Range("I21").Value = 13
$ws.Range("J21").Value = 10000
$ws.Range("K21").Value = 13
$ws.Range("L21").Value = 10000
$ws.Range("N21").Value = -10470
$ws.Range("H31").Value = 2991.72
$ws.Range("I31").Value = 1189.9
$ws.Range("K31").Value = 1189.9
$ws.Range("M31").Value = -894.9000000000001
$ws.Range("H34").Value = 2991.72
$ws.Range("I34").Value = 1189.9
$ws.Range("K34").Value = 1189.9
$ws.Range("M34").Value = -987.9000000000001
$ws.Range("H86").Value = 6000
$ws.Range("I86").Value = 6000
$ws.Range("K86").Value = 6000
$ws.Range("M86").Value = -4877
$ws.Range("H89").Value = 6000
$ws.Range("I89").Value = 6000
$ws.Range("K89").Value = 30000
$ws.Range("M89").Value = -24384
$ws.Range("H122").Value = 1396.5333
$ws.Range("J122").Value = 1662.3334
$ws.Range("L122").Value = 4987.0002
$ws.Range("N122").Value = -9887.0002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 151.84616
$ws.Range("I33").Value = 50.2
$ws.Range("J33").Value = 215.375
$ws.Range("K33").Value = 301.2
$ws.Range("L33").Value = 1292.25
$ws.Range("M33").Value = -18.20000000000005
$ws.Range("N33").Value = -1858.25
$ws.Range("H49").Value = 868.5
$ws.Range("I49").Value = 642.2
$ws.Range("K49").Value = 1926.6
$ws.Range("M49").Value = -1770.6
$ws.Range("H100").Value = 24575
$ws.Range("I100").Value = 3512.5
$ws.Range("K100").Value = 10537.5
$ws.Range("M100").Value = -9726.5
$ws.Range("H121").Value = 853.58826
$ws.Range("J121").Value = 1121
$ws.Range("L121").Value = 3363
$ws.Range("N121").Value = -5983
$ws.Range("H139").Value = 3766.913
$ws.Range("I139").Value = 3309.2666
$ws.Range("J139").Value = 4625
$ws.Range("K139").Value = 9927.799800000001
$ws.Range("L139").Value = 13875
$ws.Range("M139").Value = -4787.799800000001
$ws.Range("N139").Value = -24155

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 30408
$ws.Range("I122").Value = 30476
$ws.Range("K122").Value = 91428
$ws.Range("M122").Value = -88978

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1257
$ws.Range("I16").Value = 450
$ws.Range("K16").Value = 450
$ws.Range("M16").Value = -280
$ws.Range("H46").Value = 4612.5
$ws.Range("I46").Value = 3241.75
$ws.Range("K46").Value = 3241.75
$ws.Range("M46").Value = -3053.75
$ws.Range("H61").Value = 5213.6
$ws.Range("I61").Value = 4517
$ws.Range("K61").Value = 4517
$ws.Range("M61").Value = -4315
$ws.Range("H68").Value = 5759
$ws.Range("J68").Value = 7249.7144
$ws.Range("L68").Value = 7249.7144
$ws.Range("N68").Value = -8747.714400000001
$ws.Range("H71").Value = 5759
$ws.Range("J71").Value = 7249.7144
$ws.Range("L71").Value = 36248.572
$ws.Range("N71").Value = -43736.572
$ws.Range("H93").Value = 7513.5713
$ws.Range("I93").Value = 6795
$ws.Range("K93").Value = 6795
$ws.Range("M93").Value = -5547
$ws.Range("H100").Value = 5421.7407
$ws.Range("I100").Value = 3830.5
$ws.Range("K100").Value = 3830.5
$ws.Range("M100").Value = -3289.5
$ws.Range("H113").Value = 5213.6
$ws.Range("I113").Value = 4517
$ws.Range("K113").Value = 4517
$ws.Range("M113").Value = -2347

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 4087.8
$ws.Range("I96").Value = 2992
$ws.Range("K96").Value = 2992
$ws.Range("M96").Value = -1619
$ws.Range("H122").Value = 3736.7407
$ws.Range("I122").Value = 1844.0555
$ws.Range("J122").Value = 7522.1113
$ws.Range("K122").Value = 5532.166499999999
$ws.Range("L122").Value = 22566.3339
$ws.Range("M122").Value = -3082.166499999999
$ws.Range("N122").Value = -27466.3339
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("N128").ClearContents()
